$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 24.26000000000035
$ws.Range("H2").Value = [double]"1.820037745287142e-16"
$ws.Range("K2").Value = 47.36189377814651
$ws.Range("L2").Value = "[39.760234405061645, 54.96355315123137]"
$ws.Range("O2").Value = 1.465447624197041
$ws.Range("P2").Value = "[1.2893423303021176, 1.6415529180919641]"
$ws.Range("S2").Value = 60.37232774251761
$ws.Range("T2").Value = "[55.55797159906223, 65.18668388597298]"
$ws.Range("W2").Value = 18.60176176176203
$ws.Range("X2").Value = 17.92180180180206
$ws.Range("Y2").Value = 19.281721721722

# Row 3 updates
$ws.Range("E3").Value = 24.90000000000045
$ws.Range("H3").Value = [double]"1.820037745287142e-16"
$ws.Range("I3").Value = 0.5373576599000143
$ws.Range("K3").Value = 49.24542089406816
$ws.Range("L3").Value = "[38.03498343471918, 60.45585835341714]"
$ws.Range("M3").Value = [double]"2.220446049250313e-16"
$ws.Range("N3").Value = [double]"2.220446049250313e-16"
$ws.Range("O3").Value = 2.408868841491272
$ws.Range("P3").Value = "[2.182447749340657, 2.6352899336418876]"
$ws.Range("S3").Value = 67.55062840906882
$ws.Range("T3").Value = "[61.8626790305013, 73.23857778763634]"
$ws.Range("W3").Value = 15.35375375375404
$ws.Range("X3").Value = 14.45645645645672
$ws.Range("Y3").Value = 16.25105105105135
